$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# Change header "View" -> "Cache"
$ws.Range("F1").Value = "Cache"

# Set default value of column F (rows 2-26) to FALSE
for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 6).Value = $false
}

# Update selection to match the new state
$ws.Range("F2:F26").Select()
